# Update "want to go" counts (column F) on sheet "展览" and sheet "全部类型"
# to reflect refreshed data generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) — row -> new F value
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 6551
$wsExpo.Range("F5").Value  = 409
$wsExpo.Range("F8").Value  = 538
$wsExpo.Range("F9").Value  = 95
$wsExpo.Range("F14").Value = 1116
$wsExpo.Range("F15").Value = 3243
$wsExpo.Range("F18").Value = 1890

# Sheet "全部类型" (All types) — same events, shifted one row down from row 8 on
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 6551
$wsAll.Range("F5").Value  = 409
$wsAll.Range("F9").Value  = 538
$wsAll.Range("F10").Value = 95
$wsAll.Range("F15").Value = 1116
$wsAll.Range("F16").Value = 3243
$wsAll.Range("F19").Value = 1890
